# Apply the cryptocurrency price/volume updates for the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.341.89"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.666.53"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("E4").Value = "  +1.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("E7").Value = "  +0.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2664"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06394"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.86"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07846"
$ws.Range("D11").ClearFormats()

$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").Value = "1.667.95"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "1.894.78"
$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "0.0₅8182"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.91"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "26.355.98"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.683"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.77"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.033"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.012"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1225"
$ws.Range("D26").ClearFormats()

$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05869"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.283"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.578"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.279"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.602"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9687"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.829"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.420"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5804"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8628"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.57%  "

$ws.Range("D41").Value = "1.065.52"
$ws.Range("E41").Value = "  +3.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.817"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.67%  "

$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.23"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").Value = "1.805.39"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.24%  "

$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  -5.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.016"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.45%  "

$ws.Range("E49").Value = "  +1.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.035"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05170"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.56%  "
